$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial number. Every populated row
# (2 through 173) currently stores the same date serial 46075 (2026-02-22)
# and needs to be bumped forward by one day to 46076 (2026-02-23).
for ($row = 2; $row -le 173; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 46075) {
        $cell.Value2 = 46076
    }
}
